$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-04-29 Tuesday" "2025-04-30 Wednesday"

Replace-Text "33×74=2442" "49×51=2499"
Replace-Text "24×42=1008" "60×40=2400"
Replace-Text "62×81=5022" "96×32=3072"
Replace-Text "81×99=8019" "35×60=2100"
Replace-Text "63×43=2709" "46×89=4094"

Replace-Text "56×25=1400" "40×17=680"
Replace-Text "87×93=8091" "94×62=5828"
Replace-Text "39×72=2808" "92×51=4692"
Replace-Text "78×38=2964" "69×30=2070"
Replace-Text "39×40=1560" "73×92=6716"

Replace-Text "95×35=3325" "20×41=820"
Replace-Text "14×25=350" "40×48=1920"
Replace-Text "53×55=2915" "56×76=4256"
Replace-Text "97×36=3492" "16×62=992"
Replace-Text "57×63=3591" "12×58=696"

Replace-Text "72×76=5472" "36×53=1908"
Replace-Text "51×82=4182" "52×68=3536"
Replace-Text "85×44=3740" "50×99=4950"
Replace-Text "97×88=8536" "59×37=2183"
Replace-Text "92×96=8832" "54×56=3024"

Replace-Text "17×41=697" "95×16=1520"
Replace-Text "82×25=2050" "63×33=2079"
Replace-Text "62×25=1550" "60×18=1080"
Replace-Text "29×66=1914" "20×18=360"
Replace-Text "22×44=968" "78×53=4134"

Write-Output "Done"
